$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The scraper re-ran and re-ordered a few same-day fixtures (rows 15/16,
# 50/51 and 53/54 swap places) plus appended one brand-new match row (100).
# Columns A:E (index, pais, torneio, temporada, data_partida) are identical
# within each swapped pair, so only F:V need to be exchanged.
# ---------------------------------------------------------------------------

function Swap-MatchRows($RowA, $RowB, $ValuesA, $ValuesB) {
    # $ValuesA / $ValuesB hold the original F:V values (17 entries, F..V) of
    # $RowA / $RowB respectively; write each row's values into the other row.
    $cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $RowA).Value = $ValuesB[$i]
        $ws.Range($cols[$i] + $RowB).Value = $ValuesA[$i]
    }
}

# --- rows 15 / 16 -----------------------------------------------------------
$row15 = @("Bandirmaspor", 0, "Boluspor", 0, 1.86, "13/08/2023 16:13", 1.73, "20/08/2023 15:53", 3.66, "13/08/2023 16:13", 3.82, "20/08/2023 15:53", 4.06, "13/08/2023 16:13", 4.79, "20/08/2023 15:53", "https://www.betexplorer.com/football/turkey/1-lig/bandirmaspor-boluspor/zLm1jMIP/")
$row16 = @("Manisa FK", 0, "Keciorengucu", 0, 1.97, "13/08/2023 16:13", 1.87, "20/08/2023 15:54", 3.63, "13/08/2023 16:13", 3.58, "20/08/2023 15:59", 3.66, "13/08/2023 16:13", 4.27, "20/08/2023 15:54", "https://www.betexplorer.com/football/turkey/1-lig/manisa-fk-keciorengucu/lSTypKeg/")
Swap-MatchRows 15 16 $row15 $row16

# --- rows 50 / 51 -----------------------------------------------------------
$row50 = @("Manisa FK", 0, "Goztepe", 1, 2.07, "17/09/2023 18:13", 2.28, "23/09/2023 10:33", 3.34, "17/09/2023 18:13", 3.35, "23/09/2023 10:33", 3.45, "17/09/2023 18:13", 3.26, "23/09/2023 10:33", "https://www.betexplorer.com/football/turkey/1-lig/manisa-fk-goztepe/j7K2PnF2/")
$row51 = @("Corum", 1, "Sakaryaspor", 0, 2.01, "19/09/2023 16:13", 1.97, "23/09/2023 17:59", 3.45, "19/09/2023 16:13", 3.4, "23/09/2023 17:59", 3.73, "19/09/2023 16:13", 4.08, "23/09/2023 17:59", "https://www.betexplorer.com/football/turkey/1-lig/corum-fk-sakaryaspor/2ov835hq/")
Swap-MatchRows 50 51 $row50 $row51

# --- rows 53 / 54 -----------------------------------------------------------
$row53 = @("Sanliurfaspor", 0, "Kocaelispor", 2, 2.31, "19/09/2023 16:13", 2.76, "24/09/2023 17:57", 3.29, "19/09/2023 16:13", 3.45, "24/09/2023 17:57", 3.15, "19/09/2023 16:13", 2.55, "24/09/2023 16:28", "https://www.betexplorer.com/football/turkey/1-lig/sanliurfaspor-kocaelispor/lvbX97FF/")
$row54 = @("Adanaspor AS", 1, "Erzurumspor", 0, 2.05, "17/09/2023 18:13", 2.64, "24/09/2023 17:57", 3.43, "17/09/2023 18:13", 3.28, "24/09/2023 17:59", 3.61, "17/09/2023 18:13", 2.78, "24/09/2023 17:57", "https://www.betexplorer.com/football/turkey/1-lig/adanaspor-as-erzurumspor-fk/AmaTAm09/")
Swap-MatchRows 53 54 $row53 $row54

# ---------------------------------------------------------------------------
# New row 100: a new match appended at the bottom of the table.
# Copy formats from row 99 (A gets the bold/bordered "index" style, E gets
# the date-time number format) so no stray new style entries are created,
# then fill in the values.
# ---------------------------------------------------------------------------
$ws.Range("A99").Copy() | Out-Null
$ws.Range("A100").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("E99").Copy() | Out-Null
$ws.Range("E100").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A100").Value = 99
$ws.Range("B100").Value = "turkey"
$ws.Range("C100").Value = "1-lig"
$ws.Range("D100").Value = "2023-2024"
$ws.Range("E100").Value = 45236.75
$ws.Range("F100").Value = "Umraniyespor"
$ws.Range("G100").Value = 1
$ws.Range("H100").Value = "Tuzlaspor"
$ws.Range("I100").Value = 1
$ws.Range("J100").Value = 1.69
$ws.Range("K100").Value = "30/10/2023 18:12"
$ws.Range("L100").Value = 1.72
$ws.Range("M100").Value = "06/11/2023 17:29"
$ws.Range("N100").Value = 3.84
$ws.Range("O100").Value = "30/10/2023 18:12"
$ws.Range("P100").Value = 3.82
$ws.Range("Q100").Value = "06/11/2023 17:29"
$ws.Range("R100").Value = 4.8
$ws.Range("S100").Value = "30/10/2023 18:12"
$ws.Range("T100").Value = 4.85
$ws.Range("U100").Value = "06/11/2023 17:29"
$ws.Range("V100").Value = "https://www.betexplorer.com/football/turkey/1-lig/umraniyespor-tuzlaspor/f1dQYmlo/"
